# fix bug trong file convertJsonToExcel
# Adds two missing "Đơn sale chính" rows (HD-LUXURY 675/676), updates its
# "Tổng" summary row, adds a missing "Đơn thu nợ" row (TN 197), updates its
# "Tổng" summary row, and refreshes the dependent totals on "Lương".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Đơn sale chính" (1st sheet)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Push the current "Tổng" row (row 9) down by two rows so the two new
# data rows can be inserted above it.
$ws1.Rows.Item(9).Insert()
$ws1.Rows.Item(9).Insert()

# Make sure date-like text columns stay plain text instead of being
# auto-converted into Excel date serials.
$ws1.Columns.Item(3).NumberFormat = "@"

# New row 9
$ws1.Cells.Item(9, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(9, 2).Value = 675
$ws1.Cells.Item(9, 3).Value = "08-23-2024"
$ws1.Cells.Item(9, 4).Value = "CẦN THƠ"
$ws1.Cells.Item(9, 5).Value = "Đoàn Minh Thư"
$ws1.Cells.Item(9, 6).Value = "Cá nhân"
$ws1.Cells.Item(9, 7).Value = "Tiêm Filler"
$ws1.Cells.Item(9, 11).Value = 0
$ws1.Cells.Item(9, 12).Value = 1500000
$ws1.Cells.Item(9, 13).Value = 0.1
$ws1.Cells.Item(9, 14).Value = 0

# New row 10
$ws1.Cells.Item(10, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(10, 2).Value = 676
$ws1.Cells.Item(10, 3).Value = "08-23-2024"
$ws1.Cells.Item(10, 4).Value = "CẦN THƠ"
$ws1.Cells.Item(10, 5).Value = "Nguyễn Thị Ngọc Tâm"
$ws1.Cells.Item(10, 6).Value = "Khách cũ"
$ws1.Cells.Item(10, 7).Value = "Tiêm Filler"
$ws1.Cells.Item(10, 8).Value = 3700000
$ws1.Cells.Item(10, 11).Value = 3700000
$ws1.Cells.Item(10, 12).Value = 3700000
$ws1.Cells.Item(10, 13).Value = 0.1
$ws1.Cells.Item(10, 14).Value = 370000

# Updated "Tổng" row, now row 11
$ws1.Cells.Item(11, 1).Value = "Tổng"
$ws1.Cells.Item(11, 2).Value = 9
$ws1.Cells.Item(11, 8).Value = 28300000
$ws1.Cells.Item(11, 10).Value = 0
$ws1.Cells.Item(11, 11).Value = 28300000
$ws1.Cells.Item(11, 12).Value = 29800000
$ws1.Cells.Item(11, 13).Value = 0
$ws1.Cells.Item(11, 14).Value = 2992000

# ---------------------------------------------------------------------
# Sheet "Đơn thu nợ" (3rd sheet)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Row 2 updates
$ws3.Cells.Item(2, 15).Value = 12000000
$ws3.Cells.Item(2, 21).Value = 0

# Push the current "Tổng" row (row 3) down one row for the new data row.
$ws3.Rows.Item(3).Insert()

$ws3.Columns.Item(4).NumberFormat = "@"

# New row 3
$ws3.Cells.Item(3, 1).Value = "TN"
$ws3.Cells.Item(3, 2).Value = 197
$ws3.Cells.Item(3, 3).Value = 1000000
$ws3.Cells.Item(3, 4).Value = "08-24-2024"
$ws3.Cells.Item(3, 5).Value = "CẦN THƠ"
$ws3.Cells.Item(3, 6).Value = "HD-LUXURY-538"
$ws3.Cells.Item(3, 7).Value = "Nâng mũi"
$ws3.Cells.Item(3, 8).Value = "Ngô Xuân Nhi"
$ws3.Cells.Item(3, 9).Value = "Cá nhân"
$ws3.Cells.Item(3, 10).Value = "Lâm Hoàng Phú"
$ws3.Cells.Item(3, 11).Value = 10000000
$ws3.Cells.Item(3, 12).Value = "Đỗ Thị Huyền Trân"
$ws3.Cells.Item(3, 13).Value = 8000000
$ws3.Cells.Item(3, 14).Value = 18000000
$ws3.Cells.Item(3, 15).Value = 12000000
$ws3.Cells.Item(3, 16).Value = "Lâm Thị Mỹ Hằng"
$ws3.Cells.Item(3, 18).Value = 0
$ws3.Cells.Item(3, 19).Value = 0
$ws3.Cells.Item(3, 20).Value = 0.04
$ws3.Cells.Item(3, 21).Value = 0
$ws3.Cells.Item(3, 22).Value = 0
$ws3.Cells.Item(3, 23).Value = 0
$ws3.Cells.Item(3, 24).Value = 0
$ws3.Cells.Item(3, 25).Value = 0

# Updated "Tổng" row, now row 4
$ws3.Cells.Item(4, 1).Value = "Tổng"
$ws3.Cells.Item(4, 2).Value = 2
$ws3.Cells.Item(4, 3).Value = 2500000
$ws3.Cells.Item(4, 11).Value = 20000000
$ws3.Cells.Item(4, 13).Value = 16000000
$ws3.Cells.Item(4, 14).Value = 36000000
$ws3.Cells.Item(4, 15).Value = 24000000
$ws3.Cells.Item(4, 18).Value = 0
$ws3.Cells.Item(4, 19).Value = 0
$ws3.Cells.Item(4, 20).Value = 0
$ws3.Cells.Item(4, 21).Value = 0
$ws3.Cells.Item(4, 22).Value = 0
$ws3.Cells.Item(4, 23).Value = 0
$ws3.Cells.Item(4, 24).Value = 0
$ws3.Cells.Item(4, 25).Value = 0

# ---------------------------------------------------------------------
# Sheet "Lương" (4th sheet)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(5, 2).Value = 2992000
$ws4.Cells.Item(11, 2).Value = 0
$ws4.Cells.Item(35, 2).Value = 7880214.285714285
$ws4.Cells.Item(38, 2).Value = 7880214.285714285
